$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 64: Test 67deg, frase12, 30 sec, dataset_86, 500, 500, 50 - no accuracy/efficacy/n speech/time yet
$ws.Range("A64").Value = "67°"
$ws.Range("B64").Value = "frase12"
$ws.Range("C64").Value = "30 sec"
$ws.Range("D64").Value = "dataset_86"
$ws.Range("E64").Value = 500
$ws.Range("F64").Value = 500
$ws.Range("G64").Value = 50

# Row 65: Test 68deg, frase121, 10 sec, dataset_86, 500, 500, 50, accuracy 0.9309, efficacy 0.9308, n speech 14, time 6gg
$ws.Range("A65").Value = "68°"
$ws.Range("B65").Value = "frase121"
$ws.Range("C65").Value = "10 sec"
$ws.Range("D65").Value = "dataset_86"
$ws.Range("E65").Value = 500
$ws.Range("F65").Value = 500
$ws.Range("G65").Value = 50
$ws.Range("H65").Value = 0.9309
$ws.Range("I65").Value = 0.9308
$ws.Range("J65").Value = "14"
$ws.Range("K65").Value = "6gg"

# Row 66: Test 69deg, frase122, 10 sec, dataset_86, 500, 500, 50
$ws.Range("A66").Value = "69°"
$ws.Range("B66").Value = "frase122"
$ws.Range("C66").Value = "10 sec"
$ws.Range("D66").Value = "dataset_86"
$ws.Range("E66").Value = 500
$ws.Range("F66").Value = 500
$ws.Range("G66").Value = 50

# Row 67: Test 70deg, frase123, 10 sec, dataset_86, 500, 500, 50
$ws.Range("A67").Value = "70°"
$ws.Range("B67").Value = "frase123"
$ws.Range("C67").Value = "10 sec"
$ws.Range("D67").Value = "dataset_86"
$ws.Range("E67").Value = 500
$ws.Range("F67").Value = 500
$ws.Range("G67").Value = 50

# Update the active selection, matching the author's final cursor state
$ws.Activate()
$ws.Range("K65").Select()
